$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.629.06"
$ws.Range("E2").Value = "  +4.25%  "

$ws.Range("D3").Value = "2.557.12"
$ws.Range("E3").Value = "  +5.42%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.89"
$ws.Range("E5").Value = "  +2.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.75"
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -6.39%  "

$ws.Range("D9").Value = "2.583.74"
$ws.Range("E9").Value = "  +5.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.45"
$ws.Range("E10").Value = "  +2.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("E13").Value = "  +0.88%  "

$ws.Range("D14").Value = "3.016.06"
$ws.Range("E14").Value = "  +6.01%  "

$ws.Range("D15").Value = "59.802.70"
$ws.Range("E15").Value = "  +4.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.47"
$ws.Range("E16").Value = "  +3.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +2.52%  "

$ws.Range("D18").Value = "2.585.51"
$ws.Range("E18").Value = "  +5.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  +1.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.00"
$ws.Range("E20").Value = "  +4.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.28"
$ws.Range("E21").Value = "  +2.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.03"
$ws.Range("E22").Value = "  +1.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.57"
$ws.Range("E24").Value = "  +2.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.418"
$ws.Range("E25").Value = "  +3.73%  "

$ws.Range("D26").Value = "2.696.57"
$ws.Range("E26").Value = "  +5.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("E27").Value = "  +2.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").Value = "0.0₃0843"
$ws.Range("E29").Value = "  +5.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("E30").Value = "  +1.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.73"
$ws.Range("E32").Value = "  +3.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.11"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("E34").Value = "  +1.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.67"
$ws.Range("E35").Value = "  +6.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.94"
$ws.Range("E36").Value = "  +4.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("E37").Value = "  +4.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("E38").Value = "  +24.62%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.46"
$ws.Range("E39").Value = "  +5.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.73"
$ws.Range("E40").Value = "  +5.56%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.834"
$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "296.66"
$ws.Range("E42").Value = "  +6.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.39"
$ws.Range("E43").Value = "  +3.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0564"
$ws.Range("E44").Value = "  +5.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.615"
$ws.Range("E45").Value = "  +2.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0991"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.993"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.52"
$ws.Range("E48").Value = "  +8.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.87"
$ws.Range("E49").Value = "  +5.03%  "

$ws.Range("D50").Value = "2.026.22"
$ws.Range("E50").Value = "  +7.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0232"
$ws.Range("E51").Value = "  +1.07%  "
